$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.139.72"
$ws.Range("E2").Value = "  -3.95%  "

$ws.Range("D3").Value = "3.153.54"
$ws.Range("E3").Value = "  -9.00%  "

$ws.Range("D4").Value = "'1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.04%  "

$ws.Range("D5").Value = "'561.15"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -4.36%  "

$ws.Range("D6").Value = "'168.32"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -5.65%  "

$ws.Range("D7").Value = "'0.611"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -2.98%  "

$ws.Range("E8").Value = "  +0.07%  "

$ws.Range("D9").Value = "3.153.37"
$ws.Range("E9").Value = "  -8.94%  "

$ws.Range("E10").Value = "  -7.73%  "

$ws.Range("E11").Value = "  -6.14%  "

$ws.Range("E12").Value = "  -6.40%  "

$ws.Range("D13").Value = "3.696.60"
$ws.Range("E13").Value = "  -9.04%  "

$ws.Range("E14").Value = "  +0.63%  "

$ws.Range("D15").Value = "'27.08"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -10.22%  "

$ws.Range("D16").Value = "64.077.41"
$ws.Range("E16").Value = "  -3.71%  "

$ws.Range("E17").Value = "  -6.52%  "

$ws.Range("D18").Value = "3.154.65"
$ws.Range("E18").Value = "  -9.89%  "

$ws.Range("E19").Value = "  -4.37%  "

$ws.Range("E20").Value = "  -7.39%  "

$ws.Range("D21").Value = "'351.37"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -5.54%  "

$ws.Range("D22").Value = "'7.18"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -6.46%  "

$ws.Range("E23").Value = "  +0.38%  "

$ws.Range("D24").Value = "'68.05"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -7.39%  "

$ws.Range("D25").Value = "'0.497"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -7.31%  "

$ws.Range("D26").Value = "'0.0000115"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -9.38%  "

$ws.Range("D27").Value = "'9.56"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -4.35%  "

$ws.Range("E28").Value = "  -1.95%  "

$ws.Range("D29").Value = "'1.00"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.15%  "

$ws.Range("E30").Value = "  -0.12%  "

$ws.Range("E31").Value = "  -6.36%  "

$ws.Range("D32").Value = "'5.44"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -8.37%  "

$ws.Range("D33").Value = "'21.87"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -7.67%  "

$ws.Range("D34").Value = "'6.58"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -7.27%  "

$ws.Range("E35").Value = "  -6.57%  "

$ws.Range("E36").Value = "  -9.52%  "

$ws.Range("D37").Value = "'153.53"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -5.41%  "

$ws.Range("D38").Value = "'0.814"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -8.09%  "

$ws.Range("D39").Value = "'26.28"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -5.96%  "

$ws.Range("E40").Value = "  -7.35%  "

$ws.Range("B41").Value = "Maker"
$ws.Range("C41").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D41").Value = "2.622.40"
$ws.Range("E41").Value = "  -5.42%  "

$ws.Range("B42").Value = "dogwifhat"
$ws.Range("C42").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D42").Value = "'2.44"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -5.31%  "

$ws.Range("D43").Value = "'4.15"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -8.04%  "

$ws.Range("D44").Value = "'39.25"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -2.04%  "

$ws.Range("D45").Value = "'5.96"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -7.73%  "

$ws.Range("D46").Value = "'0.0647"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -7.00%  "

$ws.Range("D47").Value = "'23.67"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -7.32%  "

$ws.Range("D48").Value = "'320.66"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -6.10%  "

$ws.Range("D49").Value = "'0.0269"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -7.14%  "

$ws.Range("E50").Value = "  -3.33%  "

$ws.Range("D51").Value = "'1.00"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.09%  "
